# PlayerPerformance_4271.xlsx — add player-info + batting-extra scraping
#
# Target layout after edit:
#   1. Player Info        (new)
#   2. ODI Batting         (existing "ODI Batting", MATCH_CARD_LINK -> MATCH_CODE,
#                            link text -> bare match code, blank INNING_NUMBER cells
#                            fully cleared)
#   3. ODI Bowling         (existing "ODI Bowling", MATCH_CARD_LINK -> MATCH_CODE,
#                            link text -> bare match code)
#   4. ODI Batting Extra   (new)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as genuine text (never let Excel's
# automatic type-sniffing turn a numeric-looking string into a Number).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Helper: stamp the bold / bordered / centered "header" look used by the
# existing header rows (row 1 of every sheet).
# ---------------------------------------------------------------------------
function Set-HeaderStyle {
    param($range)
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# ===========================================================================
# 1. "ODI Batting" sheet — rename MATCH_CARD_LINK -> MATCH_CODE, replace the
#    howstat scorecard URLs with the bare MatchCode, drop the leftover blank
#    INNING_NUMBER cells.
# ===========================================================================
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$lastRowBatting = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBatting; $r++) {
    $linkCell = $wsBatting.Cells.Item($r, 4)
    $code = $linkCell.Text -replace ".*MatchCode=", ""
    Set-TextValue $linkCell $code

    $inningCell = $wsBatting.Cells.Item($r, 2)
    if ($inningCell.Text -eq "") {
        $inningCell.ClearContents()
    }
}

# ===========================================================================
# 2. "ODI Bowling" sheet — same MATCH_CARD_LINK -> MATCH_CODE treatment.
# ===========================================================================
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$wsBowling.Range("B1").Value = "MATCH_CODE"

$lastRowBowling = $wsBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRowBowling; $r++) {
    $linkCell = $wsBowling.Cells.Item($r, 2)
    $code = $linkCell.Text -replace ".*MatchCode=", ""
    Set-TextValue $linkCell $code
}

# ===========================================================================
# 3. New "Player Info" sheet, inserted before "ODI Batting" (i.e. first).
# ===========================================================================
$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $wsPlayerInfo.Cells.Item(1, $c).Value = $playerInfoHeaders[$c - 1]
}
Set-HeaderStyle $wsPlayerInfo.Range("A1:D1")

$playerInfoRow = @("4271", "Moeen Munir Ali", "Left Handed", "Right Arm Off Break")
for ($c = 1; $c -le $playerInfoRow.Length; $c++) {
    Set-TextValue $wsPlayerInfo.Cells.Item(2, $c) $playerInfoRow[$c - 1]
}

# ===========================================================================
# 4. New "ODI Batting Extra" sheet, appended after "ODI Bowling" (i.e. last).
# ===========================================================================
$wsExtra = $wb.Worksheets.Add($wb.Worksheets.Item($wb.Worksheets.Count), $wb.Worksheets.Item($wb.Worksheets.Count))
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $wsExtra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
Set-HeaderStyle $wsExtra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is numeric when known; NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL are
# text (and blank when unknown, same as the source data).
$extraRows = @(
    @("4469", "6", "1", "0", "14.81%", "NO"),
    @("4470", "6", "",  "",  "",       "NO"),
    @("4471", "",  "",  "",  "",       "NO"),
    @("4598", "",  "",  "",  "",       "NO"),
    @("4599", "",  "",  "",  "",       "NO"),
    @("4602", "6", "",  "",  "",       "NO"),
    @("4609", "7", "2", "0", "12.73%", "NO"),
    @("4613", "7", "2", "2", "19.11%", "NO"),
    @("4618", "6", "2", "2", "13.13%", "NO"),
    @("4619", "",  "",  "",  "",       "NO"),
    @("4620", "5", "1", "0", "2.99%",  "NO"),
    @("4622", "",  "",  "",  "",       "NO"),
    @("4663", "",  "",  "",  "",       "NO"),
    @("4666", "",  "",  "",  "",       "NO"),
    @("4698", "6", "1", "0", "4.06%",  "NO"),
    @("4699", "6", "6", "1", "14.91%", "NO"),
    @("4700", "6", "2", "4", "11.85%", "NO"),
    @("4711", "7", "1", "0", "6.60%",  "NO"),
    @("4713", "7", "3", "2", "12.88%", "NO"),
    @("4717", "7", "0", "0", "1.02%",  "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $row = $extraRows[$i]
    $r = $i + 2

    Set-TextValue $wsExtra.Cells.Item($r, 1) $row[0]

    $posCell = $wsExtra.Cells.Item($r, 2)
    if ($row[1] -eq "") {
        $posCell.ClearContents()
    } else {
        $posCell.Value = [double]$row[1]
    }

    Set-TextValue $wsExtra.Cells.Item($r, 3) $row[2]
    Set-TextValue $wsExtra.Cells.Item($r, 4) $row[3]
    Set-TextValue $wsExtra.Cells.Item($r, 5) $row[4]
    Set-TextValue $wsExtra.Cells.Item($r, 6) $row[5]
}

Write-Output "edit complete"
